$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.458.45"
$ws.Range("E2").Value = "  +4.98%  "

$ws.Range("D3").Value = "3.174.53"
$ws.Range("E3").Value = "  +1.94%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "398.84"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.43"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.549"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.58%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.618"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.98"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.82%  "

$ws.Range("E11").Value = "  +1.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0884"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.07%  "

$ws.Range("D13").Value = "3.670.67"
$ws.Range("E13").Value = "  +2.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.19"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.05"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.98%  "

$ws.Range("E16").Value = "  +7.63%  "

$ws.Range("D17").Value = "3.188.40"
$ws.Range("E17").Value = "  +2.86%  "

$ws.Range("E18").Value = "  -2.22%  "

$ws.Range("D19").Value = "54.490.45"
$ws.Range("E19").Value = "  +4.76%  "

$ws.Range("E20").Value = "  +4.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.87"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.46%  "

$ws.Range("D22").Value = "0.0₃0993"
$ws.Range("E22").Value = "  +2.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.63"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.64%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "273.26"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.23%  "

$ws.Range("E25").Value = "  +2.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.01"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.69"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.16%  "

$ws.Range("E28").Value = "  +0.72%  "

$ws.Range("E29").Value = "  -1.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.112"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.61%  "

$ws.Range("E32").Value = "  +6.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0497"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +9.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "36.86"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.15%  "

$ws.Range("E35").Value = "  +0.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.57"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.39%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.65"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +7.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.85"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +10.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.12"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +10.12%  "

$ws.Range("E41").Value = "  +1.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.290"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.27"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.59%  "

$ws.Range("E44").Value = "  +2.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.117"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.79%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.22"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.39%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.47"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.89%  "

$ws.Range("E48").Value = "  -1.31%  "

$ws.Range("D49").Value = "2.091.41"
$ws.Range("E49").Value = "  +1.83%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0341"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.98%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0504"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.47%  "
